# Adds foreign-key target annotations ("(de Table.Colonne)") throughout the
# data-dictionary section of the document, per the commit
# "rajout de l'identification des clés étrangères".
#
# Strategy: for every edit, locate a unique anchor string with
# Range.Find.Execute, then build a small sub-Range via absolute character
# offsets ($rng.Start / $rng.End) and assign to its .Text property. Setting
# .Text on a collapsed (zero-length) range positioned right after existing
# characters inherits that text's run formatting (language, bold, etc.),
# which keeps things like <w:lang w:val="fr-FR"/> / <w:lang w:val="en-US"/>
# intact without having to touch Font/Language properties explicitly.

$d = $word.ActiveDocument

function Find-Range([string]$anchor) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "anchor not found: $anchor"
    }
    return $rng
}

# ---------------------------------------------------------------------
# 1) Utilisateur — Clé étrangère : user_abonnement
#    -> ... user_abonnement (de Abonnement.Abonnement_ID)
# ---------------------------------------------------------------------
$rng = Find-Range(" : user_abonnement")
$rng.Collapse(0)
$rng.Text = " (de Abonnement.Abonnement_ID)"

# ---------------------------------------------------------------------
# 2) Artiste — Clé étrangère : artiste_label
#    -> ... artiste_label (de Label.Label_Nom)
# ---------------------------------------------------------------------
$rng = Find-Range("Clé étrangère : artiste_label")
$rng.Collapse(0)
$rng.Text = " (de Label.Label_Nom)"

# ---------------------------------------------------------------------
# 3) Label entity header — drop label_fondation from the column list
#    "(label_nom, label_pays, label_fondation)" -> "(label_nom, label_pays)"
# ---------------------------------------------------------------------
$rng = Find-Range(", label_fondation")
$rng.Text = ""

# ---------------------------------------------------------------------
# 4) Album — Clé étrangère : album_artiste_principal, album_label
#    -> album_artiste_principal (de Artiste.Artiste_ID), album_label (de Labe.Label_noml)
# ---------------------------------------------------------------------
$rng = Find-Range(" : album_artiste_principal, album_label")
$prefixLen = (" : album_artiste_principal").Length
$midPoint = $rng.Start + $prefixLen
$endPoint = $rng.End
# Insert the trailing annotation first so the midpoint offset (computed
# above, before the string grows) stays valid.
$tailRng = $d.Range($endPoint, $endPoint)
$tailRng.Text = " (de Labe.Label_noml)"
$midRng = $d.Range($midPoint, $midPoint)
$midRng.Text = " (de Artiste.Artiste_ID)"

# ---------------------------------------------------------------------
# 5) Titre — Clé étrangère : titre_album, titre_genre
#    -> titre_album (de Album.Album_ID), titre_genre (de Genre.Genre_ID)
# ---------------------------------------------------------------------
$rng = Find-Range("Clé étrangère : titre_album, titre_genre")
$prefixLen = ("Clé étrangère : titre_album").Length
$midPoint = $rng.Start + $prefixLen
$endPoint = $rng.End
$tailRng = $d.Range($endPoint, $endPoint)
$tailRng.Text = " (de Genre.Genre_ID)"
$midRng = $d.Range($midPoint, $midPoint)
$midRng.Text = " (de Album.Album_ID)"

# ---------------------------------------------------------------------
# 6) Playlist header — rename playlist_user -> playlist_user_id in the
#    column list.
# ---------------------------------------------------------------------
$rng = Find-Range("Playlist (playlist_id, playlist_user, playlist_nom, playlist_date_creation, playlist_publique)")
$innerStart = $rng.Start + ("Playlist").Length + 1
$innerRng = $d.Range($innerStart, $rng.End)
$innerRng.Text = "(playlist_id, playlist_user_id, playlist_nom, playlist_date_creation, playlist_publique)"

# ---------------------------------------------------------------------
# 7) Playlist — Clé étrangère : playlist_user
#    -> playlist_user_id (de Utilisateur.User_ID)
# ---------------------------------------------------------------------
$rng = Find-Range(" : playlist_user")
$rng.Collapse(0)
$rng.Text = "_id (de Utilisateur.User_ID)"

# ---------------------------------------------------------------------
# 8) Abonnement header — add the abonnement_id primary-key column that was
#    missing from the column list.
# ---------------------------------------------------------------------
$rng = Find-Range("Abonnement (abonnement_user_id, abonnement_type")
$insertPoint = $rng.Start + ("Abonnement (").Length
$insRng = $d.Range($insertPoint, $insertPoint)
$insRng.Text = "abonnement_id, "

# ---------------------------------------------------------------------
# 9) Abonnement — Clé primaire : abonnement_user_id -> abonnement_id
# ---------------------------------------------------------------------
$rng = Find-Range("Clé primaire : abonnement_user_id")
$prefixLen = ("Clé primaire : ").Length
$start = $rng.Start + $prefixLen
$valRng = $d.Range($start, $rng.End)
$valRng.Text = "abonnement_id"

# ---------------------------------------------------------------------
# 10) Abonnement — Clé étrangère : abonnement_user_id
#     -> abonnement_user_id (de Utilisateur.User_ID)
# ---------------------------------------------------------------------
$rng = Find-Range("Clé étrangère : abonnement_user_id")
$rng.Collapse(0)
$rng.Text = " (de Utilisateur.User_ID)"

# ---------------------------------------------------------------------
# 11 & 12) Historique — Clé étrangère : historique_user, historique_tite
#     -> historique_user (de Utilisateur.User_ID), historique_titre (de Titre.Titre_ID)
#     (also fixes the "tite" -> "titre" typo, per the diff)
# ---------------------------------------------------------------------
$rng = Find-Range("Clé étrangère : historique_user, historique_tite")
$prefixLen = ("Clé étrangère : historique_user").Length
$midPoint = $rng.Start + $prefixLen
$endPoint = $rng.End
$tailRng = $d.Range($endPoint - 4, $endPoint)
$tailRng.Text = "titre (de Titre.Titre_ID)"
$midRng = $d.Range($midPoint, $midPoint)
$midRng.Text = " (de Utilisateur.User_ID)"

# ---------------------------------------------------------------------
# 13 & 14) Notation — Clé étrangère : notation_user, notation_titre
#     -> notation_user (de Utilisateur.User_ID), notation_titre (de Titre.Titre_ID)
# ---------------------------------------------------------------------
$rng = Find-Range("Clé étrangère : notation_user, notation_titre")
$prefixLen = ("Clé étrangère : notation_user").Length
$midPoint = $rng.Start + $prefixLen
$endPoint = $rng.End
$tailRng = $d.Range($endPoint, $endPoint)
$tailRng.Text = " (de Titre.Titre_ID)"
$midRng = $d.Range($midPoint, $midPoint)
$midRng.Text = " (de Utilisateur.User_ID)"

Write-Output "done"
